# whitelabel.xlsx template update:
#  - ReplaceStrings: drop the "isRegex" column header (F1) - feature removed
#  - Rename: add header row (Path / newName) above the existing blank row
#  - plugin: untouched content (only style renumbering happens automatically)
#  - Delete: new sheet added (with a "Path" header), for a delete-by-path feature
#  - the "Delete" sheet becomes the active/selected sheet in the workbook

$wb = $excel.ActiveWorkbook

# --- ReplaceStrings sheet: remove the isRegex header cell ---
$wsReplace = $wb.Worksheets.Item("ReplaceStrings")
$wsReplace.Range("F1").Value = $null
$wsReplace.Range("F4").Select()

# --- Rename sheet: add "Path" / "newName" header row ---
$wsRename = $wb.Worksheets.Item("Rename")
$wsRename.Range("A1").Value = "Path"
$wsRename.Range("B1").Value = "newName"
$wsRename.Range("B1").Select()

# --- plugin sheet: no content change ---
$wsPlugin = $wb.Worksheets.Item("plugin")

# --- Delete sheet: brand new sheet at the end of the workbook ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsDelete = $wb.Worksheets.Add($null, $lastSheet)
$wsDelete.Name = "Delete"
$wsDelete.Range("A1").Value = "Path"
$wsDelete.Range("A1").Select()

# Make "Delete" the active/selected sheet (matches activeTab="3" in workbook.xml)
$wsDelete.Activate()
